$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.933.06"
$ws.Range("E2").Value = "  -7.47%  "

$ws.Range("D3").Value = "3.714.22"
$ws.Range("E3").Value = "  -6.61%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.46"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -5.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.01"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.91%  "

$ws.Range("D7").Value = "3.702.65"
$ws.Range("E7").Value = "  -6.73%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.628"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -8.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.710"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -9.92%  "

$ws.Range("E11").Value = "  -10.87%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "52.92"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000297"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -11.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.61"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -6.07%  "

$ws.Range("D15").Value = "4.320.48"
$ws.Range("E15").Value = "  -6.43%  "

$ws.Range("D16").Value = "3.710.01"
$ws.Range("E16").Value = "  -6.69%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.41"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -6.08%  "

$ws.Range("E18").Value = "  -2.96%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.95"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -8.95%  "

$ws.Range("E20").Value = "  -8.21%  "

$ws.Range("D21").Value = "67.753.27"
$ws.Range("E21").Value = "  -7.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "407.97"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -10.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.49"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -6.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.52"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -8.09%  "

$ws.Range("E25").Value = "  -10.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.80"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -10.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.69"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.83"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -9.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.96"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.54"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -9.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.88"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -9.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.63"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.63"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -9.28%  "

$ws.Range("E34").Value = "  -9.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "65.14"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -7.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "43.50"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -9.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "600.90"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -7.36%  "

$ws.Range("D38").Value = "0.0₃0903"
$ws.Range("E38").Value = "  -13.96%  "

$ws.Range("E39").Value = "  -7.24%  "

$ws.Range("E40").Value = "  +0.17%  "

$ws.Range("E41").Value = "  -0.08%  "

$ws.Range("E42").Value = "  -7.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.04"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -10.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.96"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -7.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0442"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -8.65%  "

$ws.Range("E46").Value = "  +0.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.28"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -12.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.72"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -13.19%  "

$ws.Range("E49").Value = "  -9.86%  "

$ws.Range("D50").Value = "2.747.02"
$ws.Range("E50").Value = "  -2.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.12"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -8.77%  "
